{"js": "// Fix the publication year in the document title line:\n// \"AUTONOMOUS INFINITY LIGHT CONTROL PREVENTION SECURITY SYSTEMS (2022) \u2013 ENSURES ...\"\n// becomes \"... (2023) \u2013 ENSURES ...\", keeping the existing bold formatting.\n\n// Make sure the edit isn't recorded as a tracked change (the target\n// revision is a clean, direct edit with no <w:ins>/<w:del> markup).\ncontext.document.changeTrackingMode = Word.ChangeTrackingMode.off;\nawait context.sync();\n\n// Locate the bold \"2022\" run that follows \"PREVENTION SECURITY SYSTEMS (\".\n// Scope the search tightly so we never touch an unrelated \"2022\"/\"2023\"\n// substring elsewhere in the document (e.g. inside a generated timestamp).\nconst anchor = context.document.body.search(\"PREVENTION SECURITY SYSTEMS (2022)\", { matchCase: true });\nanchor.load(\"text\");\nawait context.sync();\n\nif (anchor.items.length === 0) {\n  throw new Error('Could not find \"PREVENTION SECURITY SYSTEMS (2022)\" in the document body.');\n}\n\n// Within that unambiguous anchor, find just the \"2022\" substring and\n// replace it with \"2023\" in place, preserving the run's formatting\n// (bold / bold-complex-script) exactly as it was.\nconst yearRange = anchor.items[0].search(\"2022\", { matchCase: true }).getFirst();\nyearRange.insertText(\"2023\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Fix the publication year in the document title line:\n# \"AUTONOMOUS INFINITY LIGHT CONTROL PREVENTION SECURITY SYSTEMS (2022) - ENSURES ...\"\n# becomes \"... (2023) - ENSURES ...\", keeping the existing bold formatting.\n\n$d = $word.ActiveDocument\n\n# Make sure the edit isn't recorded as a tracked change (the target\n# revision is a clean, direct edit with no inserted/deleted revision marks).\n$d.TrackRevisions = $false\n\n# Locate the bold \"2022\" run using an unambiguous anchor phrase so we never\n# touch an unrelated \"2022\"/\"2023\" substring elsewhere in the document\n# (e.g. inside a generated timestamp).\n$anchor = $d.Content\n$anchorFind = $anchor.Find\n$anchorFind.ClearFormatting()\n$anchorFind.Text = \"PREVENTION SECURITY SYSTEMS (2022)\"\n$anchorFind.MatchCase = $true\n$anchorFind.MatchWholeWord = $false\n$anchorFound = $anchorFind.Execute()\n\nif (-not $anchorFound) {\n    throw \"Could not find 'PREVENTION SECURITY SYSTEMS (2022)' in the document.\"\n}\n\n# Work on a duplicate of the anchor range, narrowed down to just the \"2022\"\n# substring, and replace its text in place so the run's formatting\n# (bold / bold-complex-script) is preserved exactly as it was.\n$yearRange = $anchor.Duplicate\n$yearFind = $yearRange.Find\n$yearFind.ClearFormatting()\n$yearFind.Text = \"2022\"\n$yearFind.MatchCase = $true\n$yearFound = $yearFind.Execute()\n\nif (-not $yearFound) {\n    throw \"Could not find '2022' within the anchor range.\"\n}\n\n$yearRange.Text = \"2023\"\n"}
